$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain text (values look numeric but must keep exact formatting)
$textCells = @("D4","D5","D6","D12","D15","D21","D22","D24","D25","D27","D28","D29","D33","D34","D35","D36","D37","D40","D42","D44","D46","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "56.410.42"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "2.345.13"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "515.89"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").Value = "133.51"
$ws.Range("E6").Value = "  +3.54%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").Value = "2.341.95"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("E10").Value = "  +5.92%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +6.30%  "
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").Value = "2.759.41"
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").Value = "23.64"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "56.528.21"
$ws.Range("E16").Value = "  +2.79%  "
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").Value = "2.365.24"
$ws.Range("E18").Value = "  +3.71%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").Value = "319.49"
$ws.Range("E21").Value = "  +3.94%  "
$ws.Range("D22").Value = "6.59"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "60.58"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  +5.01%  "
$ws.Range("D27").Value = "7.70"
$ws.Range("E27").Value = "  +3.34%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "1.22"
$ws.Range("E28").Value = "  +8.74%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "170.56"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").Value = "0.0₃0733"
$ws.Range("E30").Value = "  +4.02%  "
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("D33").Value = "18.20"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").Value = "0.941"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("E38").Value = "  +3.84%  "
$ws.Range("E39").Value = "  +6.98%  "
$ws.Range("D40").Value = "37.44"
$ws.Range("E40").Value = "  +2.67%  "
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").Value = "137.84"
$ws.Range("E42").Value = "  +8.41%  "
$ws.Range("E43").Value = "  +4.17%  "
$ws.Range("D44").Value = "274.82"
$ws.Range("E44").Value = "  +8.68%  "
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").Value = "0.0929"
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "0.557"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D50").Value = "0.378"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "16.68"
$ws.Range("E51").Value = "  +0.88%  "
